# Actualización automática 2025-06-26 09:20:08
# Insert a new advisor row ("MONTESDEOCA ROBLES MARIA HILDA") in alphabetical
# position (row 33) on both the "VENTAS POR GRUPO" and "VENTA MENSUAL"
# sheets, shifting all following rows down by one, and refresh the
# "X de 52" -> "X de 53" counters on the trailing summary row of
# "VENTAS POR GRUPO" to reflect the new total of 53 advisors.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO" (columns A:R)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a blank row at position 33, pushing existing row 33 (and below)
# down to row 34 (and below).
$ws1.Rows("33:33").Insert()

# Populate the newly inserted row with the new advisor.
$ws1.Range("A33").Value = "GUERRERO FAREZ FABIAN MAURICIO"
$ws1.Range("B33").Value = "MONTESDEOCA ROBLES MARIA HILDA"
for ($col = 3; $col -le 18; $col++) {
    $ws1.Cells.Item(33, $col).Value = 0
}

# The trailing "count" summary row moved from 54 -> 55 automatically with
# the insert; update its "X de 52" labels to "X de 53" (52 -> 53 advisors).
for ($col = 3; $col -le 18; $col++) {
    $cell = $ws1.Cells.Item(55, $col)
    $cell.Value = $cell.Value2.Replace("de 52", "de 53")
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL" (columns A:G)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Insert a blank row at position 33, pushing existing row 33 (and below)
# down to row 34 (and below).
$ws2.Rows("33:33").Insert()

# Populate the newly inserted row with the new advisor.
$ws2.Range("A33").Value = "GUERRERO FAREZ FABIAN MAURICIO"
$ws2.Range("B33").Value = "MONTESDEOCA ROBLES MARIA HILDA"
for ($col = 3; $col -le 7; $col++) {
    $ws2.Cells.Item(33, $col).Value = 0
}

# Totals row (previously 54, now 55) values are unchanged by the insert.
